# Fix "Recorded By" (column G) ordering on the Session Analysis Results sheet.
# For every affected row, the recorder list (a comma-separated list such as
# "System, dnasr281@gmail.com" or "backup@backdoor.com, System, system")
# needs its first entry moved to the end of the list, e.g.
#   "System, dnasr281@gmail.com"            -> "dnasr281@gmail.com, System"
#   "backup@backdoor.com, System, system"   -> "System, system, backup@backdoor.com"
#   "backup@backdoor.com, System"           -> "System, backup@backdoor.com"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column G whose "Recorded By" value needs to be rotated
# (first item moved to the end of the comma-separated list).
$rowsToFix = @(
    2,3,5,6,8,10,11,12,13,14,15,17,18,19,20,
    21,22,24,26,28,29,31,32,34,36,37,38,39,40,41,
    43,44,45,46,47,48,50,52,54,55,57,58,60,62,63,
    64,65,66,67,69,70,71,72,73,74,76,78,80,81,82,
    83,84,85,86,90,92,93,94,96,99,101,106,107,108,109,
    110,111,112,116,118,119,120,122,125,127,132,133,134,135,136,
    137,138,142,144,145,146,148,151,153
)

$updated = 0
foreach ($row in $rowsToFix) {
    $cell = $ws.Range("G$row")
    $current = [string]$cell.Value2
    $parts = $current -split ", "
    if ($parts.Length -gt 1) {
        $rotated = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ", "
        $cell.Value = $rotated
        $updated++
    }
}

Write-Host "Updated $updated cells in column G"
